$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-15 Friday", 2) | Out-Null
$d.Content.Find.Execute("618×6=", $true, $false, $false, $false, $false, $true, 1, $false, "652×6=", 2) | Out-Null
$d.Content.Find.Execute("643×9=", $true, $false, $false, $false, $false, $true, 1, $false, "454×7=", 2) | Out-Null
$d.Content.Find.Execute("670×8=", $true, $false, $false, $false, $false, $true, 1, $false, "435×3=", 2) | Out-Null
$d.Content.Find.Execute("440×5=", $true, $false, $false, $false, $false, $true, 1, $false, "487×7=", 2) | Out-Null
$d.Content.Find.Execute("778×3=", $true, $false, $false, $false, $false, $true, 1, $false, "990×6=", 2) | Out-Null
$d.Content.Find.Execute("486×7=", $true, $false, $false, $false, $false, $true, 1, $false, "442×4=", 2) | Out-Null
$d.Content.Find.Execute("390×3=", $true, $false, $false, $false, $false, $true, 1, $false, "579×7=", 2) | Out-Null
$d.Content.Find.Execute("541×8=", $true, $false, $false, $false, $false, $true, 1, $false, "819×3=", 2) | Out-Null
$d.Content.Find.Execute("266×6=", $true, $false, $false, $false, $false, $true, 1, $false, "683×2=", 2) | Out-Null
$d.Content.Find.Execute("779×5=", $true, $false, $false, $false, $false, $true, 1, $false, "892×5=", 2) | Out-Null
$d.Content.Find.Execute("990×3=", $true, $false, $false, $false, $false, $true, 1, $false, "372×7=", 2) | Out-Null
$d.Content.Find.Execute("170×4=", $true, $false, $false, $false, $false, $true, 1, $false, "988×6=", 2) | Out-Null
$d.Content.Find.Execute("740×8=", $true, $false, $false, $false, $false, $true, 1, $false, "581×6=", 2) | Out-Null
$d.Content.Find.Execute("634×9=", $true, $false, $false, $false, $false, $true, 1, $false, "153×5=", 2) | Out-Null
$d.Content.Find.Execute("162×7=", $true, $false, $false, $false, $false, $true, 1, $false, "837×4=", 2) | Out-Null
$d.Content.Find.Execute("243×5=", $true, $false, $false, $false, $false, $true, 1, $false, "799×4=", 2) | Out-Null
$d.Content.Find.Execute("249×8=", $true, $false, $false, $false, $false, $true, 1, $false, "502×3=", 2) | Out-Null
$d.Content.Find.Execute("471×3=", $true, $false, $false, $false, $false, $true, 1, $false, "474×4=", 2) | Out-Null
$d.Content.Find.Execute("755×5=", $true, $false, $false, $false, $false, $true, 1, $false, "725×3=", 2) | Out-Null
$d.Content.Find.Execute("966×6=", $true, $false, $false, $false, $false, $true, 1, $false, "198×6=", 2) | Out-Null
$d.Content.Find.Execute("722×9=", $true, $false, $false, $false, $false, $true, 1, $false, "980×8=", 2) | Out-Null
$d.Content.Find.Execute("695×8=", $true, $false, $false, $false, $false, $true, 1, $false, "177×6=", 2) | Out-Null
$d.Content.Find.Execute("481×2=", $true, $false, $false, $false, $false, $true, 1, $false, "994×8=", 2) | Out-Null
$d.Content.Find.Execute("829×5=", $true, $false, $false, $false, $false, $true, 1, $false, "999×2=", 2) | Out-Null
$d.Content.Find.Execute("692×3=", $true, $false, $false, $false, $false, $true, 1, $false, "605×7=", 2) | Out-Null
